$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25 (columns C-I, K-L); column J stays 0.
$newValues = @{
    2 = @(9.636383089362395, 6.010696939438567, 12.62372203873266, 29.95522511061474, 38.64779909504731, 16.70132783096943, 27.24773208812358, 18.18271571215051, 9.581209994965675)
    3 = @(9.547331846453609, 6.018512708247147, 12.56804852434917, 30.04987995860623, 38.79532445480987, 16.78753589879139, 27.3525402868805, 17.57340346693814, 9.586342412800489)
    4 = @(9.493779580181203, 6.02391358239328, 12.53642413250304, 30.11986416655657, 38.90478982296467, 16.84490744038716, 27.42597975358888, 17.18994541480551, 9.591241589598825)
    5 = @(9.472259073193451, 6.02626615255789, 12.52418951357719, 30.15134521616777, 38.95410015471062, 16.86939905827029, 27.45817771589702, 17.03155992607399, 9.593677424172954)
    6 = @(9.468704435656955, 6.026665964513659, 12.52219762498908, 30.1567508497529, 38.96257070443351, 16.87353293102897, 27.46366088931607, 17.00513886306904, 9.594108423418982)
    7 = @(9.493488097085072, 6.023944695494706, 12.5362564786191, 30.12027676824757, 38.90543586308062, 16.84523324490358, 27.42640481143669, 17.18781765773217, 9.591272661368876)
    8 = @(9.60545643589057, 6.013267058169433, 12.60400026505069, 29.98538664067972, 38.69471936463211, 16.73012850509926, 27.28197561121608, 17.97469163793984, 9.582616847775876)
    9 = @(9.833080182242675, 5.997092000134165, 12.75672039351577, 29.81592852174636, 38.43328124771144, 16.53981701278573, 27.07146157320619, 19.43479264924039, 9.579512322617404)
    10 = @(10.00402756910474, 5.988096554987844, 12.88040955768009, 29.75052627163184, 38.33617289318587, 16.42183400071336, 26.96192525001791, 20.44646835737742, 9.585681079001288)
    11 = @(10.08234313936237, 5.984628325748912, 12.93902365177949, 29.73381712125017, 38.31307096993181, 16.37295756970278, 26.92205265049627, 20.89160503050753, 9.590317813588332)
    12 = @(10.11205659232406, 5.983404448177032, 12.9615432116446, 29.72937984198725, 38.3073847499272, 16.35514345028933, 26.90839759459683, 21.05786784790763, 9.592336166528021)
    13 = @(10.1056550826828, 5.983664056887558, 12.95667905262934, 29.73025121433993, 38.30847273517487, 16.35894908166667, 26.91127405469111, 21.02216444241903, 9.591889815423842)
    14 = @(10.08478666612949, 5.984525844948566, 12.94086994398621, 29.73341410620542, 38.31254162741654, 16.37147804814611, 26.92090023728559, 20.90533038535476, 9.590478608015127)
    15 = @(10.07201093232995, 5.98506535899101, 12.93122813946242, 29.73559801756964, 38.31543357309226, 16.37924296813674, 26.92698492985619, 20.83346286589264, 9.589648366996268)
    16 = @(9.99891864876226, 5.988335744383642, 12.87662505018552, 29.75188199693794, 38.33810950310286, 16.4251250856745, 26.96473241609637, 20.41706250816386, 9.585414836024949)
    17 = @(9.954204167651667, 5.990501626666534, 12.84371933935882, 29.76522321848932, 38.35744158901946, 16.45450374946072, 26.99044814650707, 20.15765281432973, 9.583286097167628)
    18 = @(9.928538741521752, 5.991806118103908, 12.82501464944077, 29.77412361943804, 38.37054274892077, 16.47185258092964, 27.00617558102412, 20.00703469638815, 9.582234054710622)
    19 = @(9.919858676303479, 5.992257891458823, 12.81872007082396, 29.77734738449228, 38.37531793160753, 16.47780390929093, 27.01166106916898, 19.95579968370265, 9.581907471544074)
    20 = @(9.958958759370089, 5.990264987959979, 12.84719935110576, 29.76367593344617, 38.3551782992987, 16.45132962868407, 26.98761365625545, 20.18541461923903, 9.583494871896951)
    21 = @(10.09091485357918, 5.984270290732724, 12.94550479122102, 29.73243368432152, 38.31126316135827, 16.36777910434118, 26.91803351570525, 20.939710798687, 9.590885995536121)
    22 = @(10.17747813970853, 5.980873775522939, 13.01163289148473, 29.72303748535031, 38.30042144497264, 16.31722333375589, 26.88097973057363, 21.41922578019476, 9.597246254628457)
    23 = @(10.13125546798799, 5.982638933030772, 12.97617179748276, 29.72703954317051, 38.30456432014721, 16.34383377337858, 26.89998170209933, 21.16457099355474, 9.593711971349734)
    24 = @(9.956809077445097, 5.990371787578816, 12.84562537255132, 29.76437162969789, 38.35619534607925, 16.45276322039028, 26.9888921924918, 20.17286810699487, 9.583399949667422)
    25 = @(9.770762320836139, 6.000959441910164, 12.71333993058136, 29.85146982316831, 38.4874937260667, 16.58748461296988, 27.12054186276638, 19.0498347364429, 9.57886729137166)
}

$cols = @("C","D","E","F","G","H","I","K","L")
foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
